$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.1831553333333333
$arr[0,3] = 0.549466
$arr[0,4] = 0.3807808465430998
$arr[0,5] = 0.3807808465430998
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 1.119424333333334
$arr[0,9] = 3.358273000000001
$arr[0,10] = 0.4918982535955602
$arr[0,11] = 0.4918982535955601
$arr[0,12] = 0.2050285369131112
$arr[0,13] = 1.845256832218
$arr[0,14] = 0.1873054334171898
$arr[0,15] = 0.1873054334171897
$ws.Range("E2:T2").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.1831553333333333
$arr[0,3] = 0.549466
$arr[0,4] = 0.3807808465430998
$arr[0,5] = 0.3807808465430998
$ws.Range("E3:J3").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.03665427988463741
$arr[0,1] = 0.03665427988463741
$arr[0,2] = 0.01527790213
$arr[0,3] = 0.13750111917
$arr[0,4] = 0.01395724772389995
$arr[0,5] = 0.01395724772389995
$ws.Range("O3:T3").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.1831553333333333
$arr[0,3] = 0.549466
$arr[0,4] = 0.3807808465430998
$arr[0,5] = 0.3807808465430998
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.4813416666666667
$arr[0,9] = 1.444025
$arr[0,10] = 0.2115115047669825
$arr[0,11] = 0.2115115047669825
$arr[0,12] = 0.08816029340555556
$arr[0,13] = 0.79344264065
$arr[0,14] = 0.08053952983877651
$arr[0,15] = 0.08053952983877649
$ws.Range("E4:T4").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.1831553333333333
$arr[0,3] = 0.549466
$arr[0,4] = 0.3807808465430998
$arr[0,5] = 0.3807808465430998
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.5915423333333333
$arr[0,9] = 1.774627
$arr[0,10] = 0.2599359617528199
$arr[0,11] = 0.25993596175282
$arr[0,12] = 0.1083441332424444
$arr[0,13] = 0.975097199182
$arr[0,14] = 0.09897863556323359
$arr[0,15] = 0.09897863556323359
$ws.Range("E5:T5").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.3553476858595785
$arr[0,1] = 0.3553476858595784
$arr[0,2] = 3
$arr[0,3] = 1
$arr[0,4] = 1.119424333333334
$arr[0,5] = 3.358273000000001
$arr[0,6] = 0.4918982535955602
$arr[0,7] = 0.4918982535955601
$arr[0,8] = 0.1913342459020001
$arr[0,9] = 1.722008213118
$arr[0,10] = 0.1747949060935504
$arr[0,11] = 0.1747949060935503
$ws.Range("I6:T6").Value = $arr

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.3553476858595785
$arr[0,1] = 0.3553476858595784
$ws.Range("I7:J7").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.03665427988463741
$arr[0,1] = 0.03665427988463741
$ws.Range("O7:P7").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.0130250135338552
$arr[0,1] = 0.0130250135338552
$ws.Range("S7:T7").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.3553476858595785
$arr[0,1] = 0.3553476858595784
$arr[0,2] = 3
$arr[0,3] = 1
$arr[0,4] = 0.4813416666666667
$arr[0,5] = 1.444025
$arr[0,6] = 0.2115115047669825
$arr[0,7] = 0.2115115047669825
$arr[0,8] = 0.08227188035000002
$arr[0,9] = 0.7404469231500002
$arr[0,10] = 0.07516012375162444
$arr[0,11] = 0.07516012375162442
$ws.Range("I8:T8").Value = $arr

$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 0.3553476858595785
$arr[0,1] = 0.3553476858595784
$arr[0,2] = 3
$arr[0,3] = 1
$arr[0,4] = 0.5915423333333333
$arr[0,5] = 1.774627
$arr[0,6] = 0.2599359617528199
$arr[0,7] = 0.25993596175282
$arr[0,8] = 0.101107598698
$arr[0,9] = 0.9099683882820001
$arr[0,10] = 0.09236764248054846
$arr[0,11] = 0.09236764248054846
$ws.Range("I9:T9").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.1002363333333333
$arr[0,1] = 0.300709
$arr[0,2] = 0.2083918342229165
$arr[0,3] = 0.2083918342229164
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 1.119424333333334
$arr[0,7] = 3.358273000000001
$arr[0,8] = 0.4918982535955602
$arr[0,9] = 0.4918982535955601
$arr[0,10] = 0.1122069906174445
$arr[0,11] = 1.009862915557
$arr[0,12] = 0.1025075793178281
$arr[0,13] = 0.1025075793178281
$ws.Range("G10:T10").Value = $arr

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 0.1002363333333333
$arr[0,1] = 0.300709
$arr[0,2] = 0.2083918342229165
$arr[0,3] = 0.2083918342229164
$ws.Range("G11:J11").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.03665427988463741
$arr[0,1] = 0.03665427988463741
$arr[0,2] = 0.008361213745
$arr[0,3] = 0.075250923705
$arr[0,4] = 0.007638452617279741
$arr[0,5] = 0.00763845261727974
$ws.Range("O11:T11").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.1002363333333333
$arr[0,1] = 0.300709
$arr[0,2] = 0.2083918342229165
$arr[0,3] = 0.2083918342229164
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 0.4813416666666667
$arr[0,7] = 1.444025
$arr[0,8] = 0.2115115047669825
$arr[0,9] = 0.2115115047669825
$arr[0,10] = 0.04824792374722223
$arr[0,11] = 0.4342313137250001
$arr[0,12] = 0.04407727043764063
$arr[0,13] = 0.04407727043764062
$ws.Range("G12:T12").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.1002363333333333
$arr[0,1] = 0.300709
$arr[0,2] = 0.2083918342229165
$arr[0,3] = 0.2083918342229164
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 0.5915423333333333
$arr[0,7] = 1.774627
$arr[0,8] = 0.2599359617528199
$arr[0,9] = 0.25993596175282
$arr[0,10] = 0.05929403450477777
$arr[0,11] = 0.533646310543
$arr[0,12] = 0.054168531850168
$arr[0,13] = 0.054168531850168
$ws.Range("G13:T13").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 1
$arr[0,1] = 0.3333333333333333
$arr[0,2] = 0.02668566666666667
$arr[0,3] = 0.080057
$arr[0,4] = 0.05547963337440523
$arr[0,5] = 0.05547963337440522
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 1.119424333333334
$arr[0,9] = 3.358273000000001
$arr[0,10] = 0.4918982535955602
$arr[0,11] = 0.4918982535955601
$arr[0,12] = 0.0298725846178889
$arr[0,13] = 0.268853261561
$arr[0,14] = 0.02729033476699189
$arr[0,15] = 0.02729033476699188
$ws.Range("E14:T14").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1
$arr[0,1] = 0.3333333333333333
$arr[0,2] = 0.02668566666666667
$arr[0,3] = 0.080057
$arr[0,4] = 0.05547963337440523
$arr[0,5] = 0.05547963337440522
$ws.Range("E15:J15").Value = $arr
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.03665427988463741
$arr[0,1] = 0.03665427988463741
$arr[0,2] = 0.002225984885
$arr[0,3] = 0.020033863965
$arr[0,4] = 0.00203356600960252
$arr[0,5] = 0.00203356600960252
$ws.Range("O15:T15").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 1
$arr[0,1] = 0.3333333333333333
$arr[0,2] = 0.02668566666666667
$arr[0,3] = 0.080057
$arr[0,4] = 0.05547963337440523
$arr[0,5] = 0.05547963337440522
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.4813416666666667
$arr[0,9] = 1.444025
$arr[0,10] = 0.2115115047669825
$arr[0,11] = 0.2115115047669825
$arr[0,12] = 0.01284492326944445
$arr[0,13] = 0.115604309425
$arr[0,14] = 0.01173458073894095
$arr[0,15] = 0.01173458073894095
$ws.Range("E16:T16").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = 1
$arr[0,1] = 0.3333333333333333
$arr[0,2] = 0.02668566666666667
$arr[0,3] = 0.080057
$arr[0,4] = 0.05547963337440523
$arr[0,5] = 0.05547963337440522
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.5915423333333333
$arr[0,9] = 1.774627
$arr[0,10] = 0.2599359617528199
$arr[0,11] = 0.25993596175282
$arr[0,12] = 0.1083441332424444
$arr[0,13] = 0.142071313739
$arr[0,14] = 0.01442115185886987
$arr[0,15] = 0.01442115185886987
$ws.Range("E17:T17").Value = $arr
